# Commit: "Changed the localize handlebars function to expose the data model
# and the calculates"
#
# Substance of the change: every `display.text` prompt on the "survey" sheet
# that referenced the bare handlebars variable `{{name}}` is rewritten to
# reference it off the exposed data model, `{{data.name}}`.
#
# (The shared-string table reshuffling / renumbering seen in the raw XML
# diff is just a side effect of Excel rewriting xl/sharedStrings.xml after
# the edit - the engine here does the same bookkeeping automatically on
# save, so we only need to touch the actual cell contents.)

$wb = $excel.ActiveWorkbook
$surveySheet = $wb.Worksheets.Item("survey")

$surveySheet.Range("G2").Value  = "What is {{data.name}}'s relationship to the household head?"
$surveySheet.Range("G3").Value  = "What is {{data.name}}'s sex?"
$surveySheet.Range("G4").Value  = "What is {{data.name}}'s age?"
$surveySheet.Range("G6").Value  = "Marital status of {{data.name}}."
$surveySheet.Range("G9").Value  = "For how many months during the last 12 months was {{data.name}} away from the household?"
$surveySheet.Range("G11").Value = "In what region was {{data.name}} born?"
$surveySheet.Range("G13").Value = "What is {{data.name}}'s main religion?"

# Reflect the editor having ended up with the cursor on G13 of the survey
# sheet, then restore the original active sheet/cell (the "initial" sheet
# stayed the selected tab throughout the real edit session).
$surveySheet.Range("G13").Select()

$initialSheet = $wb.Worksheets.Item("initial")
$initialSheet.Range("C17").Select()
